$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets list the same events (rows 2-16) and
# need their "想去人数" (F column) counts refreshed to the newer snapshot.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 218
    $ws.Range("F3").Value = 256
    $ws.Range("F5").Value = 806
    $ws.Range("F8").Value = 47
    $ws.Range("F10").Value = 104
    $ws.Range("F12").Value = 30
    $ws.Range("F14").Value = 188
    $ws.Range("F15").Value = 458
    $ws.Range("F16").Value = 40
}

# F7 ("合肥·第十三届次元之门动漫游戏博览会") ends up one apart between the
# two sheets in the refreshed snapshot: 6247 on 展览, 6248 on 全部类型.
$wb.Worksheets.Item("展览").Range("F7").Value = 6247
$wb.Worksheets.Item("全部类型").Range("F7").Value = 6248
